# "modified seven segment driver"
# Fill in the SEG_SEL / GND_CTRL_VEC / IO29-32 mapping columns (E, F, G)
# for the seven-segment driver pin table on Sheet1 (rows 16-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E16").Value = "SEG_SEL[4]"
$ws.Range("F16").Value = "GND_CTRL_VEC[0]"
$ws.Range("G16").Value = "IO8"

$ws.Range("E17").Value = "SEG_SEL[3]"
$ws.Range("F17").Value = "GND_CTRL_VEC[1]"
$ws.Range("G17").Value = "IO9"

$ws.Range("F18").Value = "GND_CTRL_VEC[2]"
$ws.Range("G18").Value = "IO10"

$ws.Range("E19").Value = "SEG_SEL[2]"
$ws.Range("F19").Value = "GND_CTRL_VEC[3]"
$ws.Range("G19").Value = "IO11"

$ws.Range("E20").Value = "SEG_SEL[6]"
$ws.Range("F20").Value = "SEG_SEL[0]"
$ws.Range("G20").Value = "IO26"

$ws.Range("E21").Value = "GND_CTRL_VEC[3]"
$ws.Range("F21").Value = "SEG_SEL[1]"
$ws.Range("G21").Value = "IO27"

$ws.Range("E22").Value = "SEG_SEL[1]"
$ws.Range("F22").Value = "SEG_SEL[2]"
$ws.Range("G22").Value = "IO28"

$ws.Range("E23").Value = "GND_CTRL_VEC[2]"
$ws.Range("F23").Value = "SEG_SEL[3]"
$ws.Range("G23").Value = "IO29"

$ws.Range("E24").Value = "GND_CTRL_VEC[1]"
$ws.Range("F24").Value = "SEG_SEL[4]"
$ws.Range("G24").Value = "IO30"

$ws.Range("E25").Value = "SEG_SEL[5]"
$ws.Range("F25").Value = "SEG_SEL[5]"
$ws.Range("G25").Value = "IO31"

$ws.Range("E26").Value = "SEG_SEL[0]"
$ws.Range("F26").Value = "SEG_SEL[6]"
$ws.Range("G26").Value = "IO32"

$ws.Range("E27").Value = "GND_CTRL_VEC[0]"

# The two new columns (E & F) now hold the longer "GND_CTRL_VEC[n]" style
# labels, so widen them to fit (matches the bestFit column resize Excel
# performs automatically when such content is typed in).
$ws.Columns("E:F").ColumnWidth = 15.66

# Leave the selection on the last edited cell, as in the authored workbook.
$ws.Range("F26").Select() | Out-Null
